$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values (becomes the old "row 3" record, with Q/R rounded) ---
$ws.Range("A2").Value = 111780813
$ws.Range("B2").Value = 56404
$ws.Range("E2").Value = 100048
$ws.Range("F2").Value = "Mindre hackspett"
$ws.Range("G2").Value = "Dryobates minor"

# I2 must stay a text cell containing "1" (not a number), with no style change.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1"
$ws.Range("I2").ClearFormats()

$ws.Range("M2").Value = "spel/sång"
$ws.Range("Q2").Value = 586624
$ws.Range("R2").Value = 6375428
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 values (becomes the old "row 2" record, with Q/R rounded) ---
$ws.Range("A3").Value = 111780818
$ws.Range("B3").Value = 56414
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = "Spillkråka"
$ws.Range("G3").Value = "Dryocopus martius"

# I3 must stay a text cell containing "2" (not a number), with no style change.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2"
$ws.Range("I3").ClearFormats()

$ws.Range("M3").Value = "lockläte, övriga läten"
$ws.Range("Q3").Value = 586624
$ws.Range("R3").Value = 6375428
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
